# ---------------------------------------------------------------------------
# lsh_data_processing/experiment_template.xlsx
# "added -r to bash + small fixes"
#
# 1) experiment_description: fix experiment 8's description text and add a
#    new experiment 9 row (recovering the old "treatment constraints
#    splitting of length of stay" description that experiment 8's row used
#    to (incorrectly) carry).
# 2) experiment_specification: fix experiment 8's inpatient_ward row
#    (transition_time_independent_splitting should be age_simple, not the
#    icu-restricted variant) and add the matching 3 state rows for the new
#    experiment 9.
# 3) run_specification: experiment 9 also belongs to run 7.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: experiment_description
# ---------------------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("experiment_description")

# Experiment 8's description was wrong - it should describe the "length of
# stay" split (3+treatment constraints), not duplicate experiment 7's text.
$wsDesc.Cells.Item(9, 2).Value = "LOS driven:3+treatment constraints splitting of length of stay in Inpatient Ward"

# Insert a new row for experiment 9, which carries the description that used
# to (incorrectly) live on experiment 8's row.
$wsDesc.Rows.Item(10).Insert()
$wsDesc.Cells.Item(10, 1).Value = 9
$wsDesc.Cells.Item(10, 2).Value = "LOS driven:7 + treatment_constraints splitting in length of stay in Inpatient Ward"
$wsDesc.Cells.Item(10, 3).Value = "base"
$wsDesc.Cells.Item(10, 4).Value = "1;4"

$wsDesc.Range("A10").Select()

# ---------------------------------------------------------------------------
# Sheet 2: experiment_specification
# ---------------------------------------------------------------------------
$wsSpec = $wb.Worksheets.Item("experiment_specification")

# Experiment 8 / inpatient_ward: transition_time_independent_splitting should
# be age_simple (not the icu-restriction variant already used by
# length_of_stay_splitting in column E).
$wsSpec.Cells.Item(24, 4).Value = "age_simple"

# Insert the 3 new state rows (home, inpatient_ward, intensive_care_unit) for
# experiment 9, right after experiment 8's block.
$wsSpec.Range("A26:A28").EntireRow.Insert()

$wsSpec.Cells.Item(26, 1).Value = 9
$wsSpec.Cells.Item(26, 2).Value = "home"
$wsSpec.Cells.Item(26, 3).Value = "length_of_stay_simple_two_weeks"
$wsSpec.Cells.Item(26, 4).Value = "age_simple"
$wsSpec.Cells.Item(26, 5).Value = "age_simple"

$wsSpec.Cells.Item(27, 1).Value = 9
$wsSpec.Cells.Item(27, 2).Value = "inpatient_ward"
$wsSpec.Cells.Item(27, 3).Value = "none"
$wsSpec.Cells.Item(27, 4).Value = "age_simple_intensive_care_unit_restriction"
$wsSpec.Cells.Item(27, 5).Value = "age_simple_intensive_care_unit_restriction"

$wsSpec.Cells.Item(28, 1).Value = 9
$wsSpec.Cells.Item(28, 2).Value = "intensive_care_unit"
$wsSpec.Cells.Item(28, 3).Value = "none"
$wsSpec.Cells.Item(28, 4).Value = "age_simple"
$wsSpec.Cells.Item(28, 5).Value = "none"

# ---------------------------------------------------------------------------
# Sheet 4: run_specification
# ---------------------------------------------------------------------------
$wsRunSpec = $wb.Worksheets.Item("run_specification")

# Run 7 now also exercises the new experiment 9.
$wsRunSpec.Cells.Item(21, 1).Value = 7
$wsRunSpec.Cells.Item(21, 2).Value = 9

$wsRunSpec.Range("A22").Select()
$wsRunSpec.Activate()

Write-Output "edit complete"
